# Apply team-member-report updates: add Andre Manz's role, and fill in
# Sprint 1/2/3 details for Aaron Riggs (back-end encryption work), matching
# the merged B:D "answer" cell layout used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- Seed new shared strings in the exact order the workbook originally
# introduced them (keeps the sharedStrings table / indices aligned). ------
$ws.Range("B58").Value = "Project Manager"
$ws.Range("B123").Value = "back end encryption implementation."
$ws.Range("B122").Value = "back end encryption completed in Python."
$ws.Range("B124").Value = "Python integration issues with the front end team."
$ws.Range("B132").Value = "Convert back end encryption software to java from python to better integrate with the front end team"
$ws.Range("B131").Value = "Started the conversion from python to java"
$ws.Range("B133").Value = "the way that java handles strings is destroying the encryption data, Java was having trouble reaching the mySQL server."
$ws.Range("B134").Value = "Was able to fix java's connection to the mySQL server"

# ===========================================================================
# Sprint 1 (rows 109-116)
# ===========================================================================

# Row 110: "Week of 2/12/2018" -- merge the whole row A:D, drop the stray
# "Week of" label that used to sit in C110.
$ws.Range("A110:B110").UnMerge()
$ws.Range("C110:D110").UnMerge()
$ws.Range("C110").Value = ""
$ws.Range("A110:D110").HorizontalAlignment = $xlCenter
$ws.Range("A110:D110").WrapText = $false
$ws.Range("A110:D110").Merge()

# Row 111: Name
$ws.Range("C111:D111").Value = ""
$ws.Range("B111:D111").HorizontalAlignment = $xlCenter
$ws.Range("B111:D111").WrapText = $false
$ws.Range("B111:D111").Merge()

# Row 112: Role
$ws.Range("C112:D112").Value = ""
$ws.Range("B112:D112").HorizontalAlignment = $xlCenter
$ws.Range("B112:D112").WrapText = $false
$ws.Range("B112:D112").Merge()

# Row 113: Tasks performed this week
$ws.Range("C113:D113").Value = ""
$ws.Range("B113:D113").HorizontalAlignment = $xlCenter
$ws.Range("B113:D113").WrapText = $true
$ws.Range("B113:D113").Merge()
$ws.Rows(113).RowHeight = 45

# Row 114: Tasks to be performed this week
$ws.Range("C114:D114").Value = ""
$ws.Range("B114:D114").HorizontalAlignment = $xlCenter
$ws.Range("B114:D114").WrapText = $true
$ws.Range("B114:D114").Merge()

# Row 115: Issues encountered
$ws.Range("C115:D115").Value = ""
$ws.Range("B115:D115").HorizontalAlignment = $xlCenter
$ws.Range("B115:D115").WrapText = $true
$ws.Range("B115:D115").Merge()
$ws.Rows(115).RowHeight = 60

# Row 116: Issues resolved
$ws.Range("C116:D116").Value = ""
$ws.Range("B116:D116").HorizontalAlignment = $xlCenter
$ws.Range("B116:D116").WrapText = $true
$ws.Range("B116:D116").Merge()

# ===========================================================================
# Sprint 2 (rows 118-125)
# ===========================================================================

# Row 119: "Week of 2/19/2018" -- merge whole row, drop stray "Week of" label
$ws.Range("A119:B119").UnMerge()
$ws.Range("C119:D119").UnMerge()
$ws.Range("C119").Value = ""
$ws.Range("A119:D119").HorizontalAlignment = $xlCenter
$ws.Range("A119:D119").WrapText = $false
$ws.Range("A119:D119").Merge()

# Row 120: Name
$ws.Range("C120:D120").Value = ""
$ws.Range("B120:D120").HorizontalAlignment = $xlCenter
$ws.Range("B120:D120").WrapText = $false
$ws.Range("B120:D120").Merge()

# Row 121: Role
$ws.Range("B121").Value = "Developer"
$ws.Range("C121:D121").Value = ""
$ws.Range("B121:D121").HorizontalAlignment = $xlCenter
$ws.Range("B121:D121").WrapText = $false
$ws.Range("B121:D121").Merge()

# Row 122: Tasks performed this week
$ws.Range("C122:D122").Value = ""
$ws.Range("B122:D122").HorizontalAlignment = $xlCenter
$ws.Range("B122:D122").WrapText = $true
$ws.Range("B122:D122").Merge()
$ws.Rows(122).RowHeight = 32.25

# Row 123: Tasks to be performed this week
$ws.Range("C123:D123").Value = ""
$ws.Range("B123:D123").HorizontalAlignment = $xlCenter
$ws.Range("B123:D123").WrapText = $true
$ws.Range("B123:D123").Merge()

# Row 124: Issues encountered
$ws.Range("C124:D124").Value = ""
$ws.Range("B124:D124").HorizontalAlignment = $xlCenter
$ws.Range("B124:D124").WrapText = $true
$ws.Range("B124:D124").Merge()
$ws.Rows(124).RowHeight = 21.75

# Row 125: Issues resolved (left blank, but centred/merged like the rest)
$ws.Range("C125:D125").Value = ""
$ws.Range("B125:D125").HorizontalAlignment = $xlCenter
$ws.Range("B125:D125").WrapText = $false
$ws.Range("B125:D125").Merge()

# ===========================================================================
# Sprint 3 (rows 127-134)
# ===========================================================================

# Row 128: "Week of 2/26/2018" -- just drop the stray "Week of" label,
# merges here (A128:B128 / C128:D128) stay as they were.
$ws.Range("C128").Value = ""

# Row 129: Name -- clear the duplicate C/D cells entirely (no merge change)
$ws.Range("C129:D129").ClearContents()

# Row 130: Role
$ws.Range("B130").Value = "Developer"
$ws.Range("C130:D130").Value = ""
$ws.Range("B130:D130").HorizontalAlignment = $xlCenter
$ws.Range("B130:D130").WrapText = $false
$ws.Range("B130:D130").Merge()

# Row 131: Tasks performed this week (centred, no wrap, matches diff s="3")
$ws.Range("C131:D131").Value = ""
$ws.Range("B131:D131").HorizontalAlignment = $xlCenter
$ws.Range("B131:D131").WrapText = $false
$ws.Range("B131:D131").Merge()

# Row 132: Tasks to be performed this week
$ws.Range("C132:D132").Value = ""
$ws.Range("B132:D132").HorizontalAlignment = $xlCenter
$ws.Range("B132:D132").WrapText = $true
$ws.Range("B132:D132").Merge()

# Row 133: Issues encountered
$ws.Range("C133:D133").Value = ""
$ws.Range("B133:D133").HorizontalAlignment = $xlCenter
$ws.Range("B133:D133").WrapText = $true
$ws.Range("B133:D133").Merge()
$ws.Rows(133).RowHeight = 43.5

# Row 134: Issues resolved
$ws.Range("C134:D134").Value = ""
$ws.Range("B134:D134").HorizontalAlignment = $xlCenter
$ws.Range("B134:D134").WrapText = $false
$ws.Range("B134:D134").Merge()

# ===========================================================================
# View state: current selection (matches the author's final cursor position)
# ===========================================================================
$ws.Range("B123:D123").Select()
